# Generate Report for Handback
#
# The localization-status report is regenerated after a handback completes:
#  - the overall/per-language "Status" flips from "Ready for handoff" to
#    "Handed back: in sync with en-US" (shared text, so it updates everywhere
#    that status string is used)
#  - each language sheet's "Latest Target File" (I2) now links to the source
#    markdown doc, "Latest Handback File" (J2) records the returned xliff file
#    name, and "Latest Handback DateTime" (K2) records when the handback
#    lands
#  - columns that now hold longer text are widened to fit

$wb = $excel.ActiveWorkbook

$ghBase = "https://github.com/OpenLocalizationTestOrg/oltest/blob/29dd8ddd3f297160131acd0b62d33d29e2feb75d/e2e/"
$docName = "b964ff0c-27b0-4326-8a2b-cb625594757d.md"

# ---------------------------------------------------------------------
# 1) Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    This string is shared by the Overview sheet (E2/F2) and both language
#    sheets' Status column (C2), so rewriting every occurrence we can find
#    flips them all consistently.
#    NOTE: use .Value2 (not .Value) for reads/writes here.
# ---------------------------------------------------------------------
$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

$overview = $wb.Worksheets.Item("Overview")
if ($overview.Range("E2").Value2 -eq $oldStatus) { $overview.Range("E2").Value2 = $newStatus }
if ($overview.Range("F2").Value2 -eq $oldStatus) { $overview.Range("F2").Value2 = $newStatus }

$zhcn = $wb.Worksheets.Item("zh-cn")
if ($zhcn.Range("C2").Value2 -eq $oldStatus) { $zhcn.Range("C2").Value2 = $newStatus }

$dede = $wb.Worksheets.Item("de-de")
if ($dede.Range("C2").Value2 -eq $oldStatus) { $dede.Range("C2").Value2 = $newStatus }

# ---------------------------------------------------------------------
# 2) zh-cn sheet: record the handback target file + handback details
# ---------------------------------------------------------------------
$zhHandoffFile = $zhcn.Range("G2").Value2

$zhcn.Hyperlinks.Add($zhcn.Range("I2"), ($ghBase + $docName), "", "", $docName)
$zhcn.Range("J2").Value2 = $zhHandoffFile
$zhcn.Range("K2").Value2 = "2016-08-13 21:16:27"

# ---------------------------------------------------------------------
# 3) de-de sheet: record the handback target file + handback details
# ---------------------------------------------------------------------
$deHandoffFile = $dede.Range("G2").Value2

$dede.Hyperlinks.Add($dede.Range("I2"), ($ghBase + $docName), "", "", $docName)
$dede.Range("J2").Value2 = $deHandoffFile
$dede.Range("K2").Value2 = "2016-08-13 21:16:37"

# ---------------------------------------------------------------------
# 4) Re-fit the columns that now hold longer text
# ---------------------------------------------------------------------
$overview.Columns.Item(5).ColumnWidth = 29.166666666666668   # E: zh-cn status
$overview.Columns.Item(6).ColumnWidth = 29.166666666666668   # F: de-de status

$zhcn.Columns.Item(3).ColumnWidth = 29.166666666666668       # C: Status
$zhcn.Columns.Item(9).ColumnWidth = 39.166666666666664       # I: Latest Target File
$zhcn.Columns.Item(10).ColumnWidth = 39.166666666666664      # J: Latest Handback File

$dede.Columns.Item(3).ColumnWidth = 29.166666666666668       # C: Status
$dede.Columns.Item(9).ColumnWidth = 39.166666666666664       # I: Latest Target File
$dede.Columns.Item(10).ColumnWidth = 39.166666666666664      # J: Latest Handback File
